# Updates the cryptos price/volume table to the latest scrape (GitHub Actions run).
# Row 32/33 also swap Fetch.AI <-> PancakeSwap (ranking order changed upstream).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '67.138.53'; Text = $true }
    @{ Cell = 'E2'; Value = '  +0.04%  '; Text = $false }
    @{ Cell = 'D3'; Value = '3.463.35'; Text = $true }
    @{ Cell = 'E3'; Value = '  -1.02%  '; Text = $false }
    @{ Cell = 'E4'; Value = '  -0.02%  '; Text = $false }
    @{ Cell = 'D5'; Value = '593.81'; Text = $true }
    @{ Cell = 'E5'; Value = '  -0.62%  '; Text = $false }
    @{ Cell = 'D6'; Value = '180.44'; Text = $true }
    @{ Cell = 'E6'; Value = '  +3.02%  '; Text = $false }
    @{ Cell = 'D7'; Value = '0.610'; Text = $true }
    @{ Cell = 'E7'; Value = '  +3.98%  '; Text = $false }
    @{ Cell = 'D8'; Value = '0.999'; Text = $true }
    @{ Cell = 'E8'; Value = '  -0.03%  '; Text = $false }
    @{ Cell = 'D9'; Value = '3.460.86'; Text = $true }
    @{ Cell = 'E9'; Value = '  -1.03%  '; Text = $false }
    @{ Cell = 'D10'; Value = '0.141'; Text = $true }
    @{ Cell = 'E10'; Value = '  +7.23%  '; Text = $false }
    @{ Cell = 'D11'; Value = '6.96'; Text = $true }
    @{ Cell = 'E11'; Value = '  -2.72%  '; Text = $false }
    @{ Cell = 'E12'; Value = '  +0.17%  '; Text = $false }
    @{ Cell = 'D13'; Value = '4.055.64'; Text = $true }
    @{ Cell = 'E13'; Value = '  -1.28%  '; Text = $false }
    @{ Cell = 'D14'; Value = '31.95'; Text = $true }
    @{ Cell = 'E14'; Value = '  +3.73%  '; Text = $false }
    @{ Cell = 'E15'; Value = '  -0.59%  '; Text = $false }
    @{ Cell = 'D16'; Value = '67.070.69'; Text = $true }
    @{ Cell = 'E16'; Value = '  -0.07%  '; Text = $false }
    @{ Cell = 'D17'; Value = '0.0000177'; Text = $true }
    @{ Cell = 'E17'; Value = '  -0.80%  '; Text = $false }
    @{ Cell = 'D18'; Value = '3.457.01'; Text = $true }
    @{ Cell = 'E18'; Value = '  -1.06%  '; Text = $false }
    @{ Cell = 'D19'; Value = '6.23'; Text = $true }
    @{ Cell = 'E19'; Value = '  -0.88%  '; Text = $false }
    @{ Cell = 'D20'; Value = '14.19'; Text = $true }
    @{ Cell = 'E20'; Value = '  -1.85%  '; Text = $false }
    @{ Cell = 'D21'; Value = '390.75'; Text = $true }
    @{ Cell = 'E21'; Value = '  -0.71%  '; Text = $false }
    @{ Cell = 'D22'; Value = '7.93'; Text = $true }
    @{ Cell = 'E22'; Value = '  -0.71%  '; Text = $false }
    @{ Cell = 'E23'; Value = '  +0.30%  '; Text = $false }
    @{ Cell = 'D24'; Value = '5.76'; Text = $true }
    @{ Cell = 'E24'; Value = '  +1.31%  '; Text = $false }
    @{ Cell = 'D25'; Value = '72.16'; Text = $true }
    @{ Cell = 'E25'; Value = '  -1.50%  '; Text = $false }
    @{ Cell = 'D26'; Value = '0.537'; Text = $true }
    @{ Cell = 'E26'; Value = '  +0.09%  '; Text = $false }
    @{ Cell = 'D27'; Value = '0.0000122'; Text = $true }
    @{ Cell = 'E27'; Value = '  +0.10%  '; Text = $false }
    @{ Cell = 'D28'; Value = '10.38'; Text = $true }
    @{ Cell = 'E28'; Value = '  +2.11%  '; Text = $false }
    @{ Cell = 'D29'; Value = '0.175'; Text = $true }
    @{ Cell = 'E29'; Value = '  -2.90%  '; Text = $false }
    @{ Cell = 'D30'; Value = '1.00'; Text = $true }
    @{ Cell = 'E30'; Value = '  +0.51%  '; Text = $false }
    @{ Cell = 'D31'; Value = '6.14'; Text = $true }
    @{ Cell = 'E31'; Value = '  +0.10%  '; Text = $false }
    @{ Cell = 'B32'; Value = 'Fetch.AI'; Text = $false }
    @{ Cell = 'C32'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; Text = $false }
    @{ Cell = 'D32'; Value = '1.40'; Text = $true }
    @{ Cell = 'E32'; Value = '  -1.48%  '; Text = $false }
    @{ Cell = 'B33'; Value = 'PancakeSwap'; Text = $false }
    @{ Cell = 'C33'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; Text = $false }
    @{ Cell = 'D33'; Value = '2.06'; Text = $true }
    @{ Cell = 'E33'; Value = '  -0.21%  '; Text = $false }
    @{ Cell = 'D34'; Value = '23.45'; Text = $true }
    @{ Cell = 'E34'; Value = '  -0.71%  '; Text = $false }
    @{ Cell = 'D35'; Value = '7.33'; Text = $true }
    @{ Cell = 'E35'; Value = '  -0.51%  '; Text = $false }
    @{ Cell = 'E36'; Value = '  +0.01%  '; Text = $false }
    @{ Cell = 'E37'; Value = '  -2.94%  '; Text = $false }
    @{ Cell = 'D38'; Value = '163.37'; Text = $true }
    @{ Cell = 'E38'; Value = '  +0.14%  '; Text = $false }
    @{ Cell = 'D39'; Value = '0.876'; Text = $true }
    @{ Cell = 'E39'; Value = '  -0.29%  '; Text = $false }
    @{ Cell = 'D40'; Value = '2.83'; Text = $true }
    @{ Cell = 'E40'; Value = '  +10.46%  '; Text = $false }
    @{ Cell = 'E41'; Value = '  -2.50%  '; Text = $false }
    @{ Cell = 'D42'; Value = '6.79'; Text = $true }
    @{ Cell = 'E42'; Value = '  -3.05%  '; Text = $false }
    @{ Cell = 'D43'; Value = '4.65'; Text = $true }
    @{ Cell = 'E43'; Value = '  -0.13%  '; Text = $false }
    @{ Cell = 'D44'; Value = '26.12'; Text = $true }
    @{ Cell = 'E44'; Value = '  +0.25%  '; Text = $false }
    @{ Cell = 'D45'; Value = '0.0721'; Text = $true }
    @{ Cell = 'E45'; Value = '  -1.27%  '; Text = $false }
    @{ Cell = 'D46'; Value = '2.742.67'; Text = $true }
    @{ Cell = 'E46'; Value = '  -1.82%  '; Text = $false }
    @{ Cell = 'D47'; Value = '26.26'; Text = $true }
    @{ Cell = 'E47'; Value = '  -4.37%  '; Text = $false }
    @{ Cell = 'D48'; Value = '41.36'; Text = $true }
    @{ Cell = 'E48'; Value = '  -2.50%  '; Text = $false }
    @{ Cell = 'D49'; Value = '0.0299'; Text = $true }
    @{ Cell = 'E49'; Value = '  -1.78%  '; Text = $false }
    @{ Cell = 'D50'; Value = '328.08'; Text = $true }
    @{ Cell = 'E50'; Value = '  -3.42%  '; Text = $false }
    @{ Cell = 'E51'; Value = '  -3.18%  '; Text = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Text) {
        # Price column holds numeric-looking text (e.g. '66.978.64'); force text so Excel
        # does not coerce it into a Number/Date and strip the formatting.
        $rng.Value = "'" + $u.Value
    } else {
        $rng.Value = $u.Value
    }
}
